$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.06440233333333333
$ws.Cells.Item(2, 8).Value = 0.193207
$ws.Cells.Item(2, 9).Value = 0.03647206354366116
$ws.Cells.Item(2, 10).Value = 0.03647206354366116
$ws.Cells.Item(2, 13).Value = 5.482938999999999
$ws.Cells.Item(2, 14).Value = 16.448817
$ws.Cells.Item(2, 15).Value = 0.1472261722051079
$ws.Cells.Item(2, 16).Value = 0.147226172205108
$ws.Cells.Item(2, 17).Value = 0.3531140651243332
$ws.Cells.Item(2, 18).Value = 3.178026586119
$ws.Cells.Item(2, 19).Value = 0.005369642307954697
$ws.Cells.Item(2, 20).Value = 0.005369642307954698
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.06440233333333333
$ws.Cells.Item(3, 8).Value = 0.193207
$ws.Cells.Item(3, 9).Value = 0.03647206354366116
$ws.Cells.Item(3, 10).Value = 0.03647206354366116
$ws.Cells.Item(3, 15).Value = 0.5993885906243068
$ws.Cells.Item(3, 16).Value = 0.5993885906243068
$ws.Cells.Item(3, 17).Value = 1.437601335784444
$ws.Cells.Item(3, 18).Value = 12.93841202206
$ws.Cells.Item(3, 19).Value = 0.02186093876459522
$ws.Cells.Item(3, 20).Value = 0.02186093876459522
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.06440233333333333
$ws.Cells.Item(4, 8).Value = 0.193207
$ws.Cells.Item(4, 9).Value = 0.03647206354366116
$ws.Cells.Item(4, 10).Value = 0.03647206354366116
$ws.Cells.Item(4, 13).Value = 9.436472999999999
$ws.Cells.Item(4, 14).Value = 28.309419
$ws.Cells.Item(4, 15).Value = 0.2533852371705853
$ws.Cells.Item(4, 16).Value = 0.2533852371705853
$ws.Cells.Item(4, 17).Value = 0.6077308796369999
$ws.Cells.Item(4, 18).Value = 5.469577916732999
$ws.Cells.Item(4, 19).Value = 0.009241482471111241
$ws.Cells.Item(4, 20).Value = 0.009241482471111241
$ws.Cells.Item(5, 9).Value = 0.8194013021867156
$ws.Cells.Item(5, 10).Value = 0.8194013021867155
$ws.Cells.Item(5, 13).Value = 5.482938999999999
$ws.Cells.Item(5, 14).Value = 16.448817
$ws.Cells.Item(5, 15).Value = 0.1472261722051079
$ws.Cells.Item(5, 16).Value = 0.147226172205108
$ws.Cells.Item(5, 17).Value = 7.933253473221999
$ws.Cells.Item(5, 18).Value = 71.39928125899799
$ws.Cells.Item(5, 19).Value = 0.1206373172208311
$ws.Cells.Item(5, 20).Value = 0.1206373172208311
$ws.Cells.Item(6, 9).Value = 0.8194013021867156
$ws.Cells.Item(6, 10).Value = 0.8194013021867155
$ws.Cells.Item(6, 15).Value = 0.5993885906243068
$ws.Cells.Item(6, 16).Value = 0.5993885906243068
$ws.Cells.Item(6, 19).Value = 0.4911397916734172
$ws.Cells.Item(6, 20).Value = 0.4911397916734171
$ws.Cells.Item(7, 9).Value = 0.8194013021867156
$ws.Cells.Item(7, 10).Value = 0.8194013021867155
$ws.Cells.Item(7, 13).Value = 9.436472999999999
$ws.Cells.Item(7, 14).Value = 28.309419
$ws.Cells.Item(7, 15).Value = 0.2533852371705853
$ws.Cells.Item(7, 16).Value = 0.2533852371705853
$ws.Cells.Item(7, 17).Value = 13.653613910754
$ws.Cells.Item(7, 18).Value = 122.882525196786
$ws.Cells.Item(7, 19).Value = 0.2076241932924674
$ws.Cells.Item(7, 20).Value = 0.2076241932924673
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.2544986666666667
$ws.Cells.Item(8, 8).Value = 0.763496
$ws.Cells.Item(8, 9).Value = 0.1441266342696234
$ws.Cells.Item(8, 10).Value = 0.1441266342696234
$ws.Cells.Item(8, 13).Value = 5.482938999999999
$ws.Cells.Item(8, 14).Value = 16.448817
$ws.Cells.Item(8, 15).Value = 0.1472261722051079
$ws.Cells.Item(8, 16).Value = 0.147226172205108
$ws.Cells.Item(8, 17).Value = 1.395400664914666
$ws.Cells.Item(8, 18).Value = 12.558605984232
$ws.Cells.Item(8, 19).Value = 0.02121921267632218
$ws.Cells.Item(8, 20).Value = 0.02121921267632218
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.2544986666666667
$ws.Cells.Item(9, 8).Value = 0.763496
$ws.Cells.Item(9, 9).Value = 0.1441266342696234
$ws.Cells.Item(9, 10).Value = 0.1441266342696234
$ws.Cells.Item(9, 15).Value = 0.5993885906243068
$ws.Cells.Item(9, 16).Value = 0.5993885906243068
$ws.Cells.Item(9, 17).Value = 5.680968440408888
$ws.Cells.Item(9, 18).Value = 51.12871596367999
$ws.Cells.Item(9, 19).Value = 0.08638786018629446
$ws.Cells.Item(9, 20).Value = 0.08638786018629446
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.2544986666666667
$ws.Cells.Item(10, 8).Value = 0.763496
$ws.Cells.Item(10, 9).Value = 0.1441266342696234
$ws.Cells.Item(10, 10).Value = 0.1441266342696234
$ws.Cells.Item(10, 13).Value = 9.436472999999999
$ws.Cells.Item(10, 14).Value = 28.309419
$ws.Cells.Item(10, 15).Value = 0.2533852371705853
$ws.Cells.Item(10, 16).Value = 0.2533852371705853
$ws.Cells.Item(10, 17).Value = 2.401569796536
$ws.Cells.Item(10, 18).Value = 21.614128168824
$ws.Cells.Item(10, 19).Value = 0.03651956140700672
$ws.Cells.Item(10, 20).Value = 0.03651956140700672
